# Apply cryptos list update (price & volume refresh), GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.496.70"
$ws.Range("E2").Value = "'  -3.33%  "
$ws.Range("D3").Value = "'3.467.88"
$ws.Range("E3").Value = "'  -2.37%  "
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("D5").Value = "'566.24"
$ws.Range("E5").Value = "'  +0.96%  "
$ws.Range("D6").Value = "'175.40"
$ws.Range("E6").Value = "'  -8.95%  "
$ws.Range("D7").Value = "'0.629"
$ws.Range("E7").Value = "'  +2.98%  "
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("E9").Value = "'  -3.18%  "
$ws.Range("E10").Value = "'  +0.93%  "
$ws.Range("D11").Value = "'53.02"
$ws.Range("E11").Value = "'  -5.75%  "
$ws.Range("E12").Value = "'  -2.46%  "
$ws.Range("D13").Value = "'9.05"
$ws.Range("E13").Value = "'  -4.94%  "
$ws.Range("D14").Value = "'4.023.76"
$ws.Range("E14").Value = "'  -2.18%  "
$ws.Range("D15").Value = "'3.464.02"
$ws.Range("E15").Value = "'  -2.43%  "
$ws.Range("E16").Value = "'  -0.56%  "
$ws.Range("D17").Value = "'18.04"
$ws.Range("E17").Value = "'  -2.16%  "
$ws.Range("D18").Value = "'65.490.29"
$ws.Range("E18").Value = "'  -3.34%  "
$ws.Range("D19").Value = "'11.93"
$ws.Range("E19").Value = "'  -0.21%  "
$ws.Range("D20").Value = "'0.984"
$ws.Range("E20").Value = "'  -1.87%  "
$ws.Range("D21").Value = "'409.85"
$ws.Range("E21").Value = "'  +0.61%  "
$ws.Range("D22").Value = "'4.09"
$ws.Range("E22").Value = "'  +2.18%  "
$ws.Range("D23").Value = "'84.58"
$ws.Range("E23").Value = "'  -1.13%  "
$ws.Range("E24").Value = "'  -1.92%  "
$ws.Range("D25").Value = "'12.55"
$ws.Range("E25").Value = "'  +4.48%  "
$ws.Range("E26").Value = "'  -5.40%  "
$ws.Range("E27").Value = "'  -6.47%  "
$ws.Range("D28").Value = "'8.84"
$ws.Range("E28").Value = "'  -0.09%  "
$ws.Range("D29").Value = "'29.86"
$ws.Range("D30").Value = "'612.63"
$ws.Range("E30").Value = "'  -11.48%  "
$ws.Range("E31").Value = "'  -8.86%  "
$ws.Range("D32").Value = "'11.48"
$ws.Range("E32").Value = "'  -2.98%  "
$ws.Range("E33").Value = "'  -4.15%  "
$ws.Range("D34").Value = "'58.57"
$ws.Range("E34").Value = "'  -3.82%  "
$ws.Range("E35").Value = "'  +7.06%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "'  -0.16%  "
$ws.Range("D37").Value = "'0.0₃0780"
$ws.Range("E37").Value = "'  -6.71%  "
$ws.Range("E38").Value = "'  -7.13%  "
$ws.Range("D39").Value = "'3.310.44"
$ws.Range("E39").Value = "'  +8.02%  "
$ws.Range("E40").Value = "'  -6.31%  "
$ws.Range("D41").Value = "'3.34"
$ws.Range("E41").Value = "'  -1.79%  "
$ws.Range("E42").Value = "'  +0.01%  "
$ws.Range("D43").Value = "'2.84"
$ws.Range("E43").Value = "'  -5.21%  "
$ws.Range("B44").Value = "'VeChain"
$ws.Range("C44").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0410"
$ws.Range("E44").Value = "'  -2.98%  "
$ws.Range("B45").Value = "'ApeXProtocol"
$ws.Range("C45").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.21"
$ws.Range("E45").Value = "'  -2.02%  "
$ws.Range("E46").Value = "'  -8.48%  "
$ws.Range("D47").Value = "'2.66"
$ws.Range("E47").Value = "'  -2.61%  "
$ws.Range("E48").Value = "'  -0.46%  "
$ws.Range("D49").Value = "'138.47"
$ws.Range("E49").Value = "'  -1.60%  "
$ws.Range("D50").Value = "'8.31"
$ws.Range("E50").Value = "'  -9.04%  "
$ws.Range("E51").Value = "'  +5.34%  "
